$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.012.18"
$ws.Range("E2").Value = "  +1.51%  "

$ws.Range("D3").Value = "1.850.02"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"

$ws.Range("E5").Value = "  +0.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "309.63"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("E7").Value = "  +2.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3676"
$ws.Range("E8").Value = "  +2.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07229"
$ws.Range("E9").Value = "  +1.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9321"
$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("E11").Value = "  +1.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07747"
$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("D13").Value = "1.831.53"
$ws.Range("E13").Value = "  -0.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.345"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.440"
$ws.Range("E15").Value = "  +1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.04"
$ws.Range("E16").Value = "  +1.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.015"
$ws.Range("E17").Value = "  +0.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008669"
$ws.Range("E18").Value = "  +1.53%  "

$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("D20").Value = "27.033.36"
$ws.Range("E20").Value = "  +1.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.51"
$ws.Range("E21").Value = "  +1.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.074"
$ws.Range("E22").Value = "  +1.17%  "

$ws.Range("E23").Value = "  +0.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.933"
$ws.Range("E24").Value = "  +1.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.20"
$ws.Range("E25").Value = "  +0.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.25"
$ws.Range("E26").Value = "  +1.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.018"
$ws.Range("E27").Value = "  +1.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.31"
$ws.Range("E28").Value = "  +0.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.976"
$ws.Range("E29").Value = "  +2.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08861"
$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.308"
$ws.Range("E31").Value = "  +4.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.181"
$ws.Range("E32").Value = "  +1.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7419"
$ws.Range("E33").Value = "  +0.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.504"
$ws.Range("E34").Value = "  +1.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.756"
$ws.Range("E35").Value = "  -2.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.115"
$ws.Range("E36").Value = "  +3.69%  "

$ws.Range("E37").Value = "  +1.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05266"
$ws.Range("E38").Value = "  +2.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.979"
$ws.Range("E39").Value = "  +1.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5274"
$ws.Range("E40").Value = "  +4.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.012"
$ws.Range("E41").Value = "  +1.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1517"
$ws.Range("E42").Value = "  +1.48%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.259"
$ws.Range("E43").Value = "  +2.24%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.64"
$ws.Range("E44").Value = "  +5.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4754"
$ws.Range("E45").Value = "  +2.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.013"
$ws.Range("E46").Value = "  +0.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.75"
$ws.Range("E47").Value = "  +3.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.610"
$ws.Range("E48").Value = "  +2.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.61"
$ws.Range("E49").Value = "  +2.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06079"
$ws.Range("E50").Value = "  +0.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8903"
$ws.Range("E51").Value = "  +4.17%  "
